# AVERAGE_1_9_matched_errors_latest.xlsx — "added ifoCAST full series evaluation"
#
# Each data row (2-24) holds one Qn target-quarters forecast errors, one value
# per forecast horizon/vintage (columns B..K). Adding the next ifoCAST vintage
# shifts every existing horizon value one column to the left (B<-C, C<-D, ...).
# Rows 2-14 gain a brand-new observation in the freed-up last column; rows 15-24
# have no further observation yet, so their row simply becomes one cell shorter.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: currently B2:K2
$vals = @(
    2.057869132359739,
    8.902115927284317,
    -8.648538155816993,
    -0.9391415561608464,
    0.6437711047260652,
    -2.322515887623191,
    -0.05848473295872768,
    -0.5096388878854663,
    -0.6500029513994404,
    -0.1237639405372229
)
for ($c = 2; $c -le 10; $c++) {
    $ws.Cells.Item(2, $c).Value = $vals[$c-1]
}
$ws.Cells.Item(2, 11).Value = -0.2655928872764969

# Row 3: currently B3:K3
$vals = @(
    6.652313087672924,
    -10.89834099542839,
    -3.188944395772239,
    -1.606031734885327,
    -4.572318727234583,
    -2.30828757257012,
    -2.759441727496859,
    -2.899805791010833,
    -2.373566780148615,
    -2.515395726887889
)
for ($c = 2; $c -le 10; $c++) {
    $ws.Cells.Item(3, $c).Value = $vals[$c-1]
}
$ws.Cells.Item(3, 11).Value = -1.872421105485825

# Row 4: currently B4:K4
$vals = @(
    -18.36749132628568,
    -10.65809472662953,
    -9.07518206574262,
    -12.04146905809188,
    -9.777437903427414,
    -10.22859205835415,
    -10.36895612186813,
    -9.842717111005909,
    -9.984546057745183,
    -9.341571436343118
)
for ($c = 2; $c -le 10; $c++) {
    $ws.Cells.Item(4, $c).Value = $vals[$c-1]
}
$ws.Cells.Item(4, 11).Value = -9.217158730013535

# Row 5: currently B5:K5
$vals = @(
    7.513167073507937,
    9.096079734394849,
    6.129792742045593,
    8.393823896710055,
    7.942669741783317,
    7.802305678269343,
    8.32854468913156,
    8.186715742392286,
    8.82969036379435,
    8.954103070123933
)
for ($c = 2; $c -le 10; $c++) {
    $ws.Cells.Item(5, $c).Value = $vals[$c-1]
}
$ws.Cells.Item(5, 11).Value = 8.556472864086315

# Row 6: currently B6:K6
$vals = @(
    0.9564081874156993,
    -2.009878804933557,
    0.2541523497309064,
    -0.1970018051958322,
    -0.3373658687098063,
    0.1888731421524112,
    0.0470441954131372,
    0.6900188168152011,
    0.8144315231447835,
    0.4168013171071663
)
for ($c = 2; $c -le 10; $c++) {
    $ws.Cells.Item(6, $c).Value = $vals[$c-1]
}
$ws.Cells.Item(6, 11).Value = 0.3406326034063205

# Row 7: currently B7:K7
$vals = @(
    -4.157449276732949,
    -1.893418122068486,
    -2.344572276995224,
    -2.484936340509199,
    -1.958697329646981,
    -2.100526276386255,
    -1.457551654984191,
    -1.333138948654609,
    -1.730769154692226,
    -1.806937868393072
)
for ($c = 2; $c -le 10; $c++) {
    $ws.Cells.Item(7, $c).Value = $vals[$c-1]
}
$ws.Cells.Item(7, 11).Value = -1.53061581027015

# Row 8: currently B8:K8
$vals = @(
    1.546611864454844,
    1.095457709528105,
    0.9550936460141312,
    1.481332656876349,
    1.339503710137075,
    1.982478331539139,
    2.106891037868721,
    1.709260831831104,
    1.633092118130258,
    1.90941417625318
)
for ($c = 2; $c -le 10; $c++) {
    $ws.Cells.Item(8, $c).Value = $vals[$c-1]
}
$ws.Cells.Item(8, 11).Value = 1.738194353941456

# Row 9: currently B9:K9
$vals = @(
    1.156631887942306,
    1.016267824428332,
    1.54250683529055,
    1.400677888551276,
    2.04365250995334,
    2.168065216282922,
    1.770435010245305,
    1.694266296544459,
    1.970588354667381,
    1.799368532355657
)
for ($c = 2; $c -le 10; $c++) {
    $ws.Cells.Item(9, $c).Value = $vals[$c-1]
}
$ws.Cells.Item(9, 11).Value = 1.952515937864398

# Row 10: currently B10:K10
$vals = @(
    -1.025188112727922,
    -0.4989491018657047,
    -0.6407780486049788,
    0.002196572797085183,
    0.1266092791266676,
    -0.2710209269109496,
    -0.3471896406117954,
    -0.07086758248887381,
    -0.2420874048005978,
    -0.08893999929185659
)
for ($c = 2; $c -le 10; $c++) {
    $ws.Cells.Item(10, $c).Value = $vals[$c-1]
}
$ws.Cells.Item(10, 11).Value = -0.3663687737149753

# Row 11: currently B11:K11
$vals = @(
    0.08364543516793629,
    -0.05818351157133772,
    0.5847911098307261,
    0.7092038161603086,
    0.3115736101226914,
    0.2354048964218456,
    0.5117269545447671,
    0.3405071322330432,
    0.4936545377417844,
    0.2162257633186657
)
for ($c = 2; $c -le 10; $c++) {
    $ws.Cells.Item(11, $c).Value = $vals[$c-1]
}
$ws.Cells.Item(11, 11).Value = 0.05323806995971928

# Row 12: currently B12:K12
$vals = @(
    -0.1538585523806955,
    0.4891160690213684,
    0.6135287753509509,
    0.2158985693133336,
    0.1397298556124878,
    0.4160519137354094,
    0.2448320914236854,
    0.3979794969324266,
    0.1205507225093079,
    -0.04243697084963852
)
for ($c = 2; $c -le 10; $c++) {
    $ws.Cells.Item(12, $c).Value = $vals[$c-1]
}
$ws.Cells.Item(12, 11).Value = -0.1665195462441563

# Row 13: currently B13:K13
$vals = @(
    0.7495351060200912,
    0.8739478123496736,
    0.4763176063120564,
    0.4001488926112106,
    0.6764709507341322,
    0.5052511284224082,
    0.6583985339311494,
    0.3809697595080307,
    0.2179820661490843,
    0.09389949075456649
)
for ($c = 2; $c -le 10; $c++) {
    $ws.Cells.Item(13, $c).Value = $vals[$c-1]
}
$ws.Cells.Item(13, 11).Value = 0.6100525277605273

# Row 14: currently B14:K14
$vals = @(
    0.03849281619118239,
    -0.3591373898464348,
    -0.4353061035472806,
    -0.158984045424359,
    -0.330203867736083,
    -0.1770564622273418,
    -0.4544852366504605,
    -0.617472930009407,
    -0.7415555054039247,
    -0.2254024683979639
)
for ($c = 2; $c -le 10; $c++) {
    $ws.Cells.Item(14, $c).Value = $vals[$c-1]
}
$ws.Cells.Item(14, 11).Value = -0.4353061035472806

# Row 15: currently B15:K15
$vals = @(
    -0.2590580299438133,
    -0.3352267436446591,
    -0.0589046855217375,
    -0.2301245078334615,
    -0.07697710232472027,
    -0.354405876747839,
    -0.5173935701067854,
    -0.6414761455013032,
    -0.1253231084953424,
    -0.3352267436446591
)
for ($c = 2; $c -le 10; $c++) {
    $ws.Cells.Item(15, $c).Value = $vals[$c-1]
}
$ws.Cells.Item(15, 11).ClearContents()

# Row 16: currently B16:J16
$vals = @(
    0.01855976243503714,
    0.2948818205579588,
    0.1236619982462347,
    0.276809403754976,
    -0.0006193706681427817,
    -0.1636070640270892,
    -0.287689639421607,
    0.2284633975843539,
    0.01855976243503714
)
for ($c = 2; $c -le 9; $c++) {
    $ws.Cells.Item(16, $c).Value = $vals[$c-1]
}
$ws.Cells.Item(16, 10).ClearContents()

# Row 17: currently B17:I17
$vals = @(
    0.1467044301255134,
    -0.0245153921862106,
    0.1286320133225306,
    -0.1487967611005881,
    -0.3117844544595345,
    -0.4358670298540523,
    0.08028600715190851,
    -0.1296176279974082
)
for ($c = 2; $c -le 8; $c++) {
    $ws.Cells.Item(17, $c).Value = $vals[$c-1]
}
$ws.Cells.Item(17, 9).ClearContents()

# Row 18: currently B18:H18
$vals = @(
    -0.1819613811903656,
    -0.02881397568162436,
    -0.3062427501047431,
    -0.4692304434636895,
    -0.5933130188582073,
    -0.07715998185224648,
    -0.2870636170015632
)
for ($c = 2; $c -le 7; $c++) {
    $ws.Cells.Item(18, $c).Value = $vals[$c-1]
}
$ws.Cells.Item(18, 8).ClearContents()

# Row 19: currently B19:G19
$vals = @(
    0.4718454808444464,
    0.1944167064213277,
    0.0314290130623813,
    -0.09265356233213651,
    0.4234994746738243,
    0.2135958395245076
)
for ($c = 2; $c -le 6; $c++) {
    $ws.Cells.Item(19, $c).Value = $vals[$c-1]
}
$ws.Cells.Item(19, 7).ClearContents()

# Row 20: currently B20:F20
$vals = @(
    -0.08594117411414147,
    -0.2489288674730878,
    -0.3730114428676057,
    0.1431415941383551,
    -0.06676204101096155
)
for ($c = 2; $c -le 5; $c++) {
    $ws.Cells.Item(20, $c).Value = $vals[$c-1]
}
$ws.Cells.Item(20, 6).ClearContents()

# Row 21: currently B21:E21
$vals = @(
    -0.07695400962807622,
    -0.201036585022594,
    0.3151164519833668,
    0.1052128168340501
)
for ($c = 2; $c -le 4; $c++) {
    $ws.Cells.Item(21, $c).Value = $vals[$c-1]
}
$ws.Cells.Item(21, 5).ClearContents()

# Row 22: currently B22:D22
$vals = @(
    -0.5068991247689255,
    0.009253912237035311,
    -0.2006497229122814
)
for ($c = 2; $c -le 3; $c++) {
    $ws.Cells.Item(22, $c).Value = $vals[$c-1]
}
$ws.Cells.Item(22, 4).ClearContents()

# Row 23: currently B23:C23
$vals = @(
    0.6215838649243215,
    0.4116802297750048
)
for ($c = 2; $c -le 2; $c++) {
    $ws.Cells.Item(23, $c).Value = $vals[$c-1]
}
$ws.Cells.Item(23, 3).ClearContents()

# Row 24: currently B24:B24
$vals = @(
    -0.2766911554241067
)
for ($c = 2; $c -le 1; $c++) {
    $ws.Cells.Item(24, $c).Value = $vals[$c-1]
}
$ws.Cells.Item(24, 2).ClearContents()
